$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("E4").Value = 0.85
$ws.Range("E6").Value = 0.05
$ws.Range("E9").Value = 0.18
$ws.Range("E11").Value = 0.67
$ws.Range("E14").Value = 0.11
$ws.Range("E16").Value = 0.18
$ws.Range("E18").Value = 0.6
$ws.Range("E21").Value = 0.08
$ws.Range("E30").Value = 0.12
$ws.Range("E31").Value = 0.32
$ws.Range("E33").Value = 0.05
$ws.Range("E34").Value = 0
$ws.Range("E38").Value = 0.29
$ws.Range("E41").Value = 0.12
$ws.Range("E60").Value = 0.08
$ws.Range("E61").Value = 0.17
$ws.Range("E64").Value = 0.09
$ws.Range("E66").Value = 0.43
$ws.Range("E70").Value = 0.03
$ws.Range("E71").Value = 0.32
$ws.Range("E83").Value = 0.01
$ws.Range("E84").Value = 0.31
$ws.Range("E93").Value = 0.33
$ws.Range("E96").Value = 0.05
$ws.Range("E108").Value = 0.34
$ws.Range("E111").Value = 0.02
$ws.Range("E120").Value = 0.11
$ws.Range("E121").Value = 0.02

$ws = $wb.Worksheets.Item(2)
$ws.Range("E8").Value = 0.05
$ws.Range("E10").Value = 0.24
$ws.Range("E11").Value = 0.29
$ws.Range("E19").Value = 0.18
$ws.Range("E21").Value = 0.37
$ws.Range("E24").Value = 0.09
$ws.Range("E26").Value = 0.16
$ws.Range("E29").Value = 0.16
$ws.Range("E31").Value = 0.06
$ws.Range("E33").Value = 0.02
$ws.Range("E36").Value = 0.26
$ws.Range("E40").Value = 0.4
$ws.Range("E41").Value = 0.01
$ws.Range("E43").Value = 0.03
$ws.Range("E45").Value = 0.11
$ws.Range("E48").Value = 0.1
$ws.Range("E51").Value = 0.68
$ws.Range("E65").Value = 0.11
$ws.Range("E66").Value = 0.47
$ws.Range("E70").Value = 0.4
$ws.Range("E71").Value = 0.19
$ws.Range("E74").Value = 0.32
$ws.Range("E75").Value = 0.08
$ws.Range("E83").Value = 0.04
$ws.Range("E84").Value = 0.3
$ws.Range("E88").Value = 0.09
$ws.Range("E90").Value = 0.75
$ws.Range("E93").Value = 0.49
$ws.Range("E96").Value = 0.09

$ws = $wb.Worksheets.Item(3)
$ws.Range("E10").Value = 0.09
$ws.Range("E11").Value = 0.54
$ws.Range("E30").Value = 0.06
$ws.Range("E31").Value = 0.32
$ws.Range("E33").Value = 0.45
$ws.Range("E35").Value = 0.24
$ws.Range("E40").Value = 0.02
$ws.Range("E41").Value = 0.22
$ws.Range("E45").Value = 0.57
$ws.Range("E46").Value = 0.09
$ws.Range("E53").Value = 0.28
$ws.Range("E56").Value = 0.05
$ws.Range("E58").Value = 0.01
$ws.Range("E61").Value = 0.11
$ws.Range("E64").Value = 0.01
$ws.Range("E66").Value = 0.15
$ws.Range("E75").Value = 0.49
$ws.Range("E76").Value = 0.03
$ws.Range("E80").Value = 0.08
$ws.Range("E81").Value = 0.18
$ws.Range("E89").Value = 0.46
$ws.Range("E90").Value = 0.02
